# Auto-generated edit script: updates Leve profit-tracking numbers
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets per the scheduled
# market-price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 3100
$ws.Range("J29").Value = 3100
$ws.Range("L29").Value = 9300
$ws.Range("N29").Value = -9862

$ws.Range("H137").Value = 1773.45
$ws.Range("I137").Value = 973.5
$ws.Range("J137").Value = 1973.4375
$ws.Range("K137").Value = 2920.5
$ws.Range("L137").Value = 5920.3125
$ws.Range("M137").Value = -370.5
$ws.Range("N137").Value = -11020.3125

$ws.Range("H138").Value = 1996.5333
$ws.Range("I138").Value = 889.2727
$ws.Range("J138").Value = 5041.5
$ws.Range("K138").Value = 2667.8181
$ws.Range("L138").Value = 15124.5
$ws.Range("M138").Value = 2472.1819
$ws.Range("N138").Value = -25404.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5882.077
$ws.Range("I32").Value = 5882.077
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 5882.077
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -5595.077
$ws.Range("N32").ClearContents()

$ws.Range("H45").Value = 3079.5833
$ws.Range("I45").Value = 631.4
$ws.Range("K45").Value = 631.4
$ws.Range("M45").Value = -254.4

$ws.Range("H55").Value = 32821.2
$ws.Range("J55").Value = 32821.2
$ws.Range("L55").Value = 32821.2
$ws.Range("N55").Value = -33451.2

$ws.Range("H97").Value = 2130.889
$ws.Range("I97").Value = 1620
$ws.Range("J97").Value = 3459.2
$ws.Range("K97").Value = 1620
$ws.Range("L97").Value = 3459.2
$ws.Range("M97").Value = -1124
$ws.Range("N97").Value = -4451.2

$ws.Range("H122").Value = 1826.3334
$ws.Range("I122").Value = 1822.7142
$ws.Range("K122").Value = 5468.142599999999
$ws.Range("M122").Value = -3018.142599999999

$ws.Range("H123").Value = 60429
$ws.Range("J123").Value = 60429
$ws.Range("L123").Value = 60429
$ws.Range("N123").Value = -70229

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1501.4
$ws.Range("I94").Value = 1228.2727
$ws.Range("K94").Value = 1228.2727
$ws.Range("M94").Value = -777.2727

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1522.3334
$ws.Range("I31").Value = 1761.9166
$ws.Range("J31").Value = 1043.1666
$ws.Range("K31").Value = 1761.9166
$ws.Range("L31").Value = 1043.1666
$ws.Range("M31").Value = -1466.9166
$ws.Range("N31").Value = -1633.1666

$ws.Range("H34").Value = 1522.3334
$ws.Range("I34").Value = 1761.9166
$ws.Range("J34").Value = 1043.1666
$ws.Range("K34").Value = 1761.9166
$ws.Range("L34").Value = 1043.1666
$ws.Range("M34").Value = -1559.9166
$ws.Range("N34").Value = -1447.1666

$ws.Range("H86").Value = 3789.6
$ws.Range("I86").Value = 3657.6667
$ws.Range("K86").Value = 3657.6667
$ws.Range("M86").Value = -2534.6667

$ws.Range("H89").Value = 3789.6
$ws.Range("I89").Value = 3657.6667
$ws.Range("K89").Value = 18288.3335
$ws.Range("M89").Value = -12672.3335

$ws.Range("H106").Value = 39750
$ws.Range("J106").Value = 39750
$ws.Range("L106").Value = 39750
$ws.Range("N106").Value = -42274

$ws.Range("H122").Value = 1561
$ws.Range("I122").Value = 1429.4
$ws.Range("K122").Value = 4288.200000000001
$ws.Range("M122").Value = -1838.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 186.35715
$ws.Range("I11").Value = 131.53847
$ws.Range("J11").Value = 899
$ws.Range("K11").Value = 394.61541
$ws.Range("L11").Value = 2697
$ws.Range("M11").Value = -254.61541
$ws.Range("N11").Value = -2977

$ws.Range("H36").Value = 599
$ws.Range("I36").Value = 599
$ws.Range("K36").Value = 1797
$ws.Range("M36").Value = -1628

$ws.Range("H92").Value = 341.75
$ws.Range("J92").Value = 490
$ws.Range("L92").Value = 1470
$ws.Range("N92").Value = -3966

$ws.Range("H121").Value = 1957534.1
$ws.Range("I121").Value = 3017.125
$ws.Range("J121").Value = 3260545.5
$ws.Range("K121").Value = 9051.375
$ws.Range("L121").Value = 9781636.5
$ws.Range("M121").Value = -7741.375
$ws.Range("N121").Value = -9784256.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 16000
$ws.Range("J18").Value = 16000
$ws.Range("L18").Value = 16000
$ws.Range("N18").Value = -16586

$ws.Range("H80").Value = 11447
$ws.Range("I80").Value = 1946.5
$ws.Range("J80").Value = 20947.5
$ws.Range("K80").Value = 1946.5
$ws.Range("L80").Value = 20947.5
$ws.Range("M80").Value = -948.5
$ws.Range("N80").Value = -22943.5

$ws.Range("H83").Value = 11447
$ws.Range("I83").Value = 1946.5
$ws.Range("J83").Value = 20947.5
$ws.Range("K83").Value = 9732.5
$ws.Range("L83").Value = 104737.5
$ws.Range("M83").Value = -4740.5
$ws.Range("N83").Value = -114721.5

$ws.Range("H113").Value = 2540.7222
$ws.Range("J113").Value = 3499.8
$ws.Range("L113").Value = 3499.8
$ws.Range("N113").Value = -7839.8

$ws.Range("H122").Value = 3467.353
$ws.Range("I122").Value = 3947.4546
$ws.Range("J122").Value = 2587.1667
$ws.Range("K122").Value = 11842.3638
$ws.Range("L122").Value = 7761.500100000001
$ws.Range("M122").Value = -9392.363799999999
$ws.Range("N122").Value = -12661.5001

$ws.Range("H132").Value = 3241.3333
$ws.Range("I132").Value = 3362
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 10086
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -7556
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3973.6
$ws.Range("I7").Value = 3806.3333
$ws.Range("J7").Value = 4224.5
$ws.Range("K7").Value = 3806.3333
$ws.Range("L7").Value = 4224.5
$ws.Range("M7").Value = -3694.3333
$ws.Range("N7").Value = -4448.5

$ws.Range("H22").Value = 4021.5386
$ws.Range("I22").Value = 3255.7144
$ws.Range("K22").Value = 3255.7144
$ws.Range("M22").Value = -2960.7144

$ws.Range("H27").Value = 4021.5386
$ws.Range("I27").Value = 3255.7144
$ws.Range("K27").Value = 3255.7144
$ws.Range("M27").Value = -3148.7144

$ws.Range("H40").Value = 5885.316
$ws.Range("I40").Value = 4623.364
$ws.Range("K40").Value = 4623.364
$ws.Range("M40").Value = -4487.364

$ws.Range("H82").Value = 3104
$ws.Range("I82").Value = 2717.6667
$ws.Range("K82").Value = 2717.6667
$ws.Range("M82").Value = -2356.6667

$ws.Range("H85").Value = 3104
$ws.Range("I85").Value = 2717.6667
$ws.Range("K85").Value = 2717.6667
$ws.Range("M85").Value = -1469.6667

$ws.Range("H122").Value = 3616.4211
$ws.Range("I122").Value = 3358.7
$ws.Range("J122").Value = 3902.7778
$ws.Range("K122").Value = 10076.1
$ws.Range("L122").Value = 11708.3334
$ws.Range("M122").Value = -7626.099999999999
$ws.Range("N122").Value = -16608.3334

$ws.Range("H126").Value = 3973.6
$ws.Range("I126").Value = 3806.3333
$ws.Range("J126").Value = 4224.5
$ws.Range("K126").Value = 11418.9999
$ws.Range("L126").Value = 12673.5
$ws.Range("M126").Value = -8948.999899999999
$ws.Range("N126").Value = -17613.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 151.5
$ws.Range("I8").Value = 151.5
$ws.Range("K8").Value = 151.5
$ws.Range("M8").Value = -11.5

$ws.Range("H100").Value = 2836
$ws.Range("I100").Value = 502
$ws.Range("K100").Value = 1004
$ws.Range("M100").Value = -463

$ws.Range("H122").Value = 2852.75
$ws.Range("I122").Value = 2484.0527
$ws.Range("K122").Value = 7452.158100000001
$ws.Range("M122").Value = -5002.158100000001

$ws.Range("H136").Value = 676.5294
$ws.Range("I136").Value = 574
$ws.Range("K136").Value = 1722
$ws.Range("M136").Value = 828
